# Insert a new data row at row 460 (shifting existing rows 460:538 down to 461:539)
# and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(460).Insert()

$ws.Range("A460").Value = 6
$ws.Range("B460").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C460").Value = 'Metropolitana'
$ws.Range("D460").Value = 44984
$ws.Range("E460").Value = 13
$ws.Range("F460").Value = 100112032
$ws.Range("G460").Value = 'Zapallo italiano'
$ws.Range("H460").Value = 'Sin especificar'
$ws.Range("I460").Value = 'Primera'
$ws.Range("J460").Value = 620
$ws.Range("K460").Value = 7500
$ws.Range("L460").Value = 8000
$ws.Range("M460").Value = 7782
$ws.Range("N460").Value = '$/caja 50 unidades'
$ws.Range("O460").Value = 'Región de Arica y Parinacota'
$ws.Range("P460").Value = 156
$ws.Range("Q460").Value = 50
$ws.Range("R460").Value = 'Hortaliza'
